# Rename the sole worksheet from "Property1" to "DataNode" — this workbook
# unifies the DataNode / DataTable / Entity naming convention, and the sheet
# tab name needs to track that rename.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "DataNode"

# Restore the author's last-saved cursor position on the sheet.
$ws.Range("C36").Select()
